$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("B2").Value = 0.00000001157407407407407
$ws.Range("C2").Value = 0.000002233796296296296

# Add new row 3
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.000002164351851851852

# Add new row 4
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = 0.000000001493055555555556
$ws.Range("C4").Value = 0.000002094537037037037

# Apply the same number format style used in B2:C2 to the new cells
$ws.Range("B3:C4").NumberFormat = $ws.Range("B2").NumberFormat
